$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Feuil2")

function Format-CardCell($row, $col) {
    $c = $ws2.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.HorizontalAlignment = -4131
    $c.VerticalAlignment = -4108
}

function Set-CardText($row, $col, $value) {
    # Value entered while the cell is already text-formatted -> stored as text
    Format-CardCell $row $col
    $ws2.Cells.Item($row, $col).Value = $value
}

function Set-CardValue($row, $col, $value) {
    # Value entered first (keeps its native type), format applied afterwards
    $ws2.Cells.Item($row, $col).Value = $value
    Format-CardCell $row $col
}

# --- New row 2 : "TEST" card ---
Set-CardValue 2 1 "TEST"
Set-CardValue 2 2 320
Set-CardValue 2 3 "Room"
Set-CardValue 2 4 "Black"
Set-CardValue 2 5 4
Set-CardValue 2 9 "Black"
Set-CardValue 2 10 470
Set-CardValue 2 19 "S"

# --- New row 3 : "TEST123" card ---
Set-CardValue 3 1 "TEST123"
Set-CardValue 3 2 470
Set-CardValue 3 3 "Action"
Set-CardValue 3 4 "Gray"
Set-CardValue 3 7 "[SPONTANEOUS] [RESTRICT 1] : Restore 1 activator, but exhaust 2 cards you own."
Set-CardValue 3 9 "Gray"
Set-CardText  3 19 "3"

# --- Column width for the Effect 1 column (G), now wide enough for long effect text ---
$ws2.Columns.Item(7).ColumnWidth = 67.83072916666667

# --- Page setup ---
$ws2.PageSetup.Orientation = 1

# --- Final selection, matching the author's end state ---
$ws2.Range("G5").Select()
